$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 30/12/2025 14:25:02"
$ws1.Range("A3").Value = "Total filas: 285"

$ws1.Cells.Item(274, 2).Value = "14:24:51"
$ws1.Cells.Item(274, 3).Value = "14:34"
$ws1.Cells.Item(274, 4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(274, 5).Value = 10
$ws1.Cells.Item(274, 6).Value = "LP1912"
$ws1.Cells.Item(274, 7).Value = "30/12/2025"

$ws1.Cells.Item(275, 2).Value = "14:24:51"
$ws1.Cells.Item(275, 3).Value = "14:44"
$ws1.Cells.Item(275, 4).Value = "14_ABASTO"
$ws1.Cells.Item(275, 5).Value = 20
$ws1.Cells.Item(275, 6).Value = "LP1912"
$ws1.Cells.Item(275, 7).Value = "30/12/2025"

$ws1.Cells.Item(276, 2).Value = "14:24:51"
$ws1.Cells.Item(276, 3).Value = "14:56"
$ws1.Cells.Item(276, 4).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(276, 5).Value = 32
$ws1.Cells.Item(276, 6).Value = "LP1912"
$ws1.Cells.Item(276, 7).Value = "30/12/2025"

$ws1.Cells.Item(277, 2).Value = "14:24:51"
$ws1.Cells.Item(277, 3).Value = "14:58"
$ws1.Cells.Item(277, 4).Value = "215B_EL PATO"
$ws1.Cells.Item(277, 5).Value = 34
$ws1.Cells.Item(277, 6).Value = "LP1912"
$ws1.Cells.Item(277, 7).Value = "30/12/2025"

$ws1.Cells.Item(278, 2).Value = "14:24:51"
$ws1.Cells.Item(278, 3).Value = "15:00"
$ws1.Cells.Item(278, 4).Value = "81_EL PELIGRO"
$ws1.Cells.Item(278, 5).Value = 36
$ws1.Cells.Item(278, 6).Value = "LP1912"
$ws1.Cells.Item(278, 7).Value = "30/12/2025"

$ws1.Cells.Item(279, 2).Value = "14:24:51"
$ws1.Cells.Item(279, 3).Value = "15:05"
$ws1.Cells.Item(279, 4).Value = "10_OLMOS"
$ws1.Cells.Item(279, 5).Value = 41
$ws1.Cells.Item(279, 6).Value = "LP1912"
$ws1.Cells.Item(279, 7).Value = "30/12/2025"

$ws1.Cells.Item(280, 2).Value = "14:24:51"
$ws1.Cells.Item(280, 3).Value = "15:20"
$ws1.Cells.Item(280, 4).Value = "15_ABASTO"
$ws1.Cells.Item(280, 5).Value = 56
$ws1.Cells.Item(280, 6).Value = "LP1912"
$ws1.Cells.Item(280, 7).Value = "30/12/2025"

$ws1.Cells.Item(281, 2).Value = "14:24:51"
$ws1.Cells.Item(281, 3).Value = "15:21"
$ws1.Cells.Item(281, 4).Value = "26_HERNANDEZ"
$ws1.Cells.Item(281, 5).Value = 57
$ws1.Cells.Item(281, 6).Value = "LP1912"
$ws1.Cells.Item(281, 7).Value = "30/12/2025"

$ws1.Cells.Item(282, 2).Value = "14:24:51"
$ws1.Cells.Item(282, 3).Value = "15:37"
$ws1.Cells.Item(282, 4).Value = "10_OLMOS"
$ws1.Cells.Item(282, 5).Value = 73
$ws1.Cells.Item(282, 6).Value = "LP1912"
$ws1.Cells.Item(282, 7).Value = "30/12/2025"

$ws1.Cells.Item(283, 2).Value = "14:24:51"
$ws1.Cells.Item(283, 3).Value = "15:45"
$ws1.Cells.Item(283, 4).Value = "14_ABASTO"
$ws1.Cells.Item(283, 5).Value = 81
$ws1.Cells.Item(283, 6).Value = "LP1912"
$ws1.Cells.Item(283, 7).Value = "30/12/2025"

$ws1.Cells.Item(284, 2).Value = "14:24:51"
$ws1.Cells.Item(284, 3).Value = "15:45"
$ws1.Cells.Item(284, 4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(284, 5).Value = 81
$ws1.Cells.Item(284, 6).Value = "LP1912"
$ws1.Cells.Item(284, 7).Value = "30/12/2025"

$ws1.Cells.Item(285, 2).Value = "14:24:51"
$ws1.Cells.Item(285, 3).Value = "16:00"
$ws1.Cells.Item(285, 4).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(285, 5).Value = 96
$ws1.Cells.Item(285, 6).Value = "LP1912"
$ws1.Cells.Item(285, 7).Value = "30/12/2025"

$ws1.Cells.Item(286, 2).Value = "14:24:51"
$ws1.Cells.Item(286, 3).Value = "16:01"
$ws1.Cells.Item(286, 4).Value = "10_OLMOS"
$ws1.Cells.Item(286, 5).Value = 97
$ws1.Cells.Item(286, 6).Value = "LP1912"
$ws1.Cells.Item(286, 7).Value = "30/12/2025"

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 30/12/2025 14:25:02"
$ws2.Range("A3").Value = "Total filas: 21"

$ws2.Cells.Item(22, 2).Value = "30/12/2025"
$ws2.Cells.Item(22, 3).Value = "14:24:51"
$ws2.Cells.Item(22, 4).Value = "14:58"
$ws2.Cells.Item(22, 5).Value = "215B_EL PATO"
$ws2.Cells.Item(22, 6).Value = 34
$ws2.Cells.Item(22, 7).Value = "LP1912"

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 30/12/2025 14:25:02"
$ws3.Range("A3").Value = "Total filas: 42"

$ws3.Cells.Item(42, 2).Value = "30/12/2025"
$ws3.Cells.Item(42, 3).Value = "14:24:56"
$ws3.Cells.Item(42, 4).Value = "14:52"
$ws3.Cells.Item(42, 5).Value = "215D_LA PLATA"
$ws3.Cells.Item(42, 6).Value = 28
$ws3.Cells.Item(42, 7).Value = "L6203"

$ws3.Cells.Item(43, 2).Value = "30/12/2025"
$ws3.Cells.Item(43, 3).Value = "14:25:01"
$ws3.Cells.Item(43, 4).Value = "15:34"
$ws3.Cells.Item(43, 5).Value = "215A_LA PLATA"
$ws3.Cells.Item(43, 6).Value = 69
$ws3.Cells.Item(43, 7).Value = "L6173"
